# Update gh-pages data output: adjust "想去人数" (F) and "最低票价" (G)
# values on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G4").Value = 40
    $ws.Range("G5").Value = 65
    $ws.Range("F6").Value = 620
    $ws.Range("F10").Value = 397
    $ws.Range("F17").Value = 1064
    $ws.Range("F18").Value = 1432
    $ws.Range("F19").Value = 314
    $ws.Range("F26").Value = 243
    $ws.Range("F27").Value = 272
    $ws.Range("F29").Value = 1657
    $ws.Range("F33").Value = 603
    $ws.Range("F35").Value = 3862
    $ws.Range("F37").Value = 449
    $ws.Range("F39").Value = 974
    $ws.Range("F40").Value = 81
}
